$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking scrape refresh: Price (D) and Volume(1h) (E) columns are
# plain text cells (t="inlineStr" in the source). Force column D to Text
# format first so numeric-looking prices ("219.48", "0.268", ...) are not
# auto-coerced into numbers by the COM Value setter, matching the original
# text storage. Style is reset back to Normal afterwards so no formatting
# change is left behind - only the literal text content changes.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '28.321.97'
$ws.Range('E2').Value = '  +4.09%  '
$ws.Range('D3').Value = '1.732.33'
$ws.Range('E3').Value = '  +2.69%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = '219.48'
$ws.Range('E5').Value = '  +1.59%  '
$ws.Range('E6').Value = '  +0.46%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('E8').Value = '  +6.10%  '
$ws.Range('D9').Value = '0.268'
$ws.Range('E9').Value = '  +2.73%  '
$ws.Range('E10').Value = '  +1.44%  '
$ws.Range('D11').Value = '0.0895'
$ws.Range('E11').Value = '  +0.48%  '
$ws.Range('D12').Value = '1.976.38'
$ws.Range('E12').Value = '  +2.70%  '
$ws.Range('D13').Value = '1.737.57'
$ws.Range('E13').Value = '  +2.95%  '
$ws.Range('D14').Value = '4.27'
$ws.Range('E14').Value = '  +1.85%  '
$ws.Range('E15').Value = '  +2.04%  '
$ws.Range('D16').Value = '67.85'
$ws.Range('E16').Value = '  +0.91%  '
$ws.Range('D17').Value = '28.314.07'
$ws.Range('E17').Value = '  +4.06%  '
$ws.Range('D18').Value = '243.75'
$ws.Range('E18').Value = '  +2.05%  '
$ws.Range('D20').Value = '7.97'
$ws.Range('E20').Value = '  -2.60%  '
$ws.Range('E21').Value = '  -0.17%  '
$ws.Range('E22').Value = '  +1.98%  '
$ws.Range('E23').Value = '  +1.48%  '
$ws.Range('E24').Value = '  -0.42%  '
$ws.Range('D25').Value = '149.40'
$ws.Range('E25').Value = '  +0.70%  '
$ws.Range('E26').Value = '  +3.20%  '
$ws.Range('D27').Value = '16.66'
$ws.Range('E27').Value = '  +0.94%  '
$ws.Range('E28').Value = '  +0.86%  '
$ws.Range('E29').Value = '  -0.17%  '
$ws.Range('D30').Value = '0.0518'
$ws.Range('E30').Value = '  +3.09%  '
$ws.Range('E31').Value = '  +2.38%  '
$ws.Range('D33').Value = '3.29'
$ws.Range('E33').Value = '  +1.39%  '
$ws.Range('D34').Value = '1.489.15'
$ws.Range('E34').Value = '  -5.33%  '
$ws.Range('E35').Value = '  -1.70%  '
$ws.Range('D36').Value = '0.978'
$ws.Range('E36').Value = '  +2.26%  '
$ws.Range('E37').Value = '  -0.02%  '
$ws.Range('D38').Value = '2.41'
$ws.Range('E38').Value = '  +0.72%  '
$ws.Range('E39').Value = '  +1.01%  '
$ws.Range('E40').Value = '  +0.48%  '
$ws.Range('D41').Value = '70.31'
$ws.Range('E41').Value = '  +0.80%  '
$ws.Range('E42').Value = '  -0.12%  '
$ws.Range('D43').Value = '5.66'
$ws.Range('E43').Value = '  +0.43%  '
$ws.Range('E44').Value = '  +1.85%  '
$ws.Range('D45').Value = '1.880.12'
$ws.Range('E45').Value = '  +2.43%  '
$ws.Range('E46').Value = '  +1.10%  '
$ws.Range('E47').Value = '  +7.95%  '
$ws.Range('D48').Value = '0.0₆0114'
$ws.Range('E48').Value = '  +5.62%  '
$ws.Range('D49').Value = '90.88'
$ws.Range('E49').Value = '  -0.52%  '
$ws.Range('D50').Value = '8.26'
$ws.Range('E50').Value = '  +0.36%  '
$ws.Range('D51').Value = '0.105'
$ws.Range('E51').Value = '  -0.61%  '

$ws.Range("D2:D51").Style = "Normal"
